# edit.ps1
# Applies the "Updated cryptos list" data refresh to the cryptocurrency
# tracking worksheet. This mirrors a scheduled GitHub Actions job that
# re-scrapes current prices/volume percentages from coinranking.com and
# rewrites the Price (D) and Volume(1h) (E) columns; two rows (47/48)
# also swap which coin (Decentraland vs EnergySwap) occupies that rank,
# so their Coin name / Link / Price / Volume cells are updated too.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price values as plain text (e.g. "27.612.76" or "1.002").
# Many of these strings are numeric-looking, and Excel's COM layer will
# silently coerce them into real numbers (losing formatting / precision)
# unless the cell is explicitly marked as Text first.
$ws.Range("D2").Value = '27.612.76'
$ws.Range("E2").Value = '  -1.25%  '
$ws.Range("D3").Value = '1.841.93'
$ws.Range("E3").Value = '  -1.01%  '
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '314.46'
$ws.Range("E5").Value = '  -1.10%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.002'
$ws.Range("E6").Value = '  +0.19%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4221'
$ws.Range("E7").Value = '  -3.43%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3637'
$ws.Range("E8").Value = '  -1.88%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '45.05'
$ws.Range("E9").Value = '  -0.32%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07291'
$ws.Range("E10").Value = '  -2.94%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.8886'
$ws.Range("E11").Value = '  -5.31%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '20.68'
$ws.Range("E12").Value = '  -2.82%  '
$ws.Range("D13").Value = '1.836.44'
$ws.Range("E13").Value = '  -1.38%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.561'
$ws.Range("E14").Value = '  -2.56%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.330'
$ws.Range("E15").Value = '  -2.12%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.06892'
$ws.Range("E16").Value = '  +0.89%  '
$ws.Range("E17").Value = '  +0.15%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '79.17'
$ws.Range("E18").Value = '  -3.11%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000008892'
$ws.Range("E19").Value = '  -1.46%  '
$ws.Range("E20").Value = '  +0.02%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '15.43'
$ws.Range("E21").Value = '  -3.29%  '
$ws.Range("D22").Value = '27.606.91'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.982'
$ws.Range("E23").Value = '  -2.60%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '10.57'
$ws.Range("E24").Value = '  -4.54%  '
$ws.Range("D25").Value = '2.041.34'
$ws.Range("E25").Value = '  -2.41%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.915'
$ws.Range("E26").Value = '  -4.63%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '153.80'
$ws.Range("E27").Value = '  -0.55%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '19.09'
$ws.Range("E28").Value = '  +3.59%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '122.71'
$ws.Range("E29").Value = '  +8.01%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.280'
$ws.Range("E30").Value = '  -2.89%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.891'
$ws.Range("E31").Value = '  +8.47%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.08934'
$ws.Range("E32").Value = '  -0.68%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.7677'
$ws.Range("E33").Value = '  -5.66%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.571'
$ws.Range("E34").Value = '  -5.26%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.929'
$ws.Range("E35").Value = '  -0.10%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.095'
$ws.Range("E36").Value = '  -6.93%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.000'
$ws.Range("E37").Value = '  +0.04%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.102'
$ws.Range("E38").Value = '  -1.11%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05372'
$ws.Range("E39").Value = '  -2.31%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01938'
$ws.Range("E40").Value = '  -1.99%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.828'
$ws.Range("E41").Value = '  -2.58%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '6.869'
$ws.Range("E42").Value = '  -2.57%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.5089'
$ws.Range("E43").Value = '  -3.33%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.1653'
$ws.Range("E44").Value = '  -2.29%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '8.262'
$ws.Range("E45").Value = '  -6.11%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.06599'
$ws.Range("E46").Value = '  -2.74%  '
$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '10.43'
$ws.Range("E47").Value = '  -2.09%  '
$ws.Range("B48").Value = 'Decentraland'
$ws.Range("C48").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.4729'
$ws.Range("E48").Value = '  -3.67%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '104.18'
$ws.Range("E49").Value = '  -2.00%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.001'
$ws.Range("E50").Value = '  +0.18%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.633'
$ws.Range("E51").Value = '  -2.93%  '
